$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6141542792320251
$ws.Range("B1").Value = 2.157603740692139
$ws.Range("C1").Value = 2.065754652023315
$ws.Range("D1").Value = 1.823393583297729
$ws.Range("E1").Value = 0.9601841568946838
